$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from an existing header cell (e.g. H1) so I1/J1 match
# the bold/centered/bordered header style used by the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Column I: constant value of 1 for rows 2-8
# Column J: mirrors column H values for rows 2-8
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value = $hVal
}
